$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Barrio Larreynaga"
$ws.Range("A3").Value = "Ciudad Jardín"
$ws.Range("A4").Value = "Barrio San Judas"
$ws.Range("A5").Value = "Barrio La Primavera"
$ws.Range("A6").Value = "Bello Horizonte"
